$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1) - new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the formatting used by the rest of the header row (bold, centered,
# top-aligned, thin box border) so the new headers look consistent with
# the existing ones (e.g. "Salary" in AA1 / "Unnamed: 27" in AB1).
$headerRange = $ws.Range("AC1:AE1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Fill team record values (Wins=68, Losses=46, Ties=0) for all data rows (2-33)
for ($r = 2; $r -le 33; $r++) {
    $ws.Cells.Item($r, 29).Value = 68
    $ws.Cells.Item($r, 30).Value = 46
    $ws.Cells.Item($r, 31).Value = 0
}
